$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (match style of the existing header row, e.g. H1)
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Fill data rows 2-40: I = 1 (constant), J = same value as H
for ($r = 2; $r -le 40; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
